$d = $word.ActiveDocument

# --- 1. Merge "Suggest the next food item to eat (eg. ... lately)" into a
#        single run (removes the proofErr-wrapped "eg." split). ---
$text1 = "Suggest the next food item to eat (eg. Some vegetable if user did not eat vegetables lately)"
$d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# --- 2. Merge "Show trends in eg. Calorie intake, vitamins, sugars" into a
#        single run. ---
$text2 = "Show trends in eg. Calorie intake, vitamins, sugars"
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# --- 3. Merge "Notify the user through Android notifications when the user
#        should eat (eg. 8am for breakfast, 5PM for dinner)" into a single
#        run. ---
$text3 = "Notify the user through Android notifications when the user should eat (eg. 8am for breakfast, 5PM for dinner)"
$d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# --- 4. Insert new bullet ideas under "Book Review" section, right after
#        "The reviews would be viewable by other users" and before the
#        "What Can I Cook?" heading. ---
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*The reviews would be viewable by other users*") {
        $anchor = $p
        break
    }
}

$newItems = @(
    "User favourites a book they liked",
    "User chooses to dislike a book they read",
    "User can share a list of their favourites/dislikes.",
    "User can choose to not be recommended books by a certain author",
    "Separate books into categories",
    "User can choose their primary language",
    "If the book has no copyright, provide a link to it or have an in-app download"
)

$current = $anchor
foreach ($item in $newItems) {
    $current.Range.InsertParagraphAfter()
    $current = $current.Next()
    $current.Range.InsertAfter($item)
}

# --- 5. Merge "Show the end user a list of recipes they can do after they
#        input the ingredients they have in the fridge" into a single run
#        (removes the gramStart/gramEnd-wrapped "ingredients" split). ---
$text5 = "Show the end user a list of recipes they can do after they input the ingredients they have in the fridge"
$d.Content.Find.Execute($text5, $true, $false, $false, $false, $false, $true, 1, $false, $text5, 2) | Out-Null

Write-Output "done"
